$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) column values - force text format to preserve exact string representation
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.564.60'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.958.34'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.79'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.71'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.956.65'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.502'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.31'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.143'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000233'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.40'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.449.11'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.411.73'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.959.29'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '441.24'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.48'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.673'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.67'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.03'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.95'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.13'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.37'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0875'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.992'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.60'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.92'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.282'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.23'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.715.42'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '135.15'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0342'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '364.10'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.05'

# Update Volume(1h) (E) column values
$ws.Range("E2").Value = '  +2.11%  '
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("E10").Value = '  +4.38%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("E20").Value = '  +0.90%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("E22").Value = '  -1.13%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -3.77%  '
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("E31").Value = '  -6.63%  '
$ws.Range("E32").Value = '  -1.43%  '
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -2.35%  '
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("E39").Value = '  +1.56%  '
$ws.Range("E40").Value = '  -3.64%  '
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("E42").Value = '  -5.36%  '
$ws.Range("E43").Value = '  -2.07%  '
$ws.Range("E44").Value = '  -7.09%  '
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("E47").Value = '  -1.89%  '
$ws.Range("E48").Value = '  -3.60%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("E51").Value = '  -4.25%  '
